# RBA v2.5 - Atualizacao da Tela
# Replace placeholder tokens (TERE/TRE/Tre/tre family) with QWER/Qwer/Qewr/qwer family,
# matching each occurrence individually (order matters, since not all occurrences of a
# given case-variant map to the same replacement).

$d = $word.ActiveDocument

# --- Word constants used below ---
# wdFindStop    = 0   (no wrap)
# wdReplaceOne  = 1
# wdCollapseEnd = 0

# 1) Main document body: the single bold "TERE" right after "A " becomes "QWER".
$bodyRng = $d.Content
$bodyRng.Find.Execute("TERE", $true, $true, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null

# 2) Header: several "TRE"/"TERE"/"Tre"/"tre" tokens must be replaced in document order,
#    each with its own (sometimes different) replacement text.
$hdr = $d.Sections.Item(1).Headers.Item(1)
$hdrEnd = $hdr.Range.End

$rng = $hdr.Range.Duplicate
$rng.Collapse(1)

$headerReplacements = @(
    @{ Find = "TRE";  Replace = "QWER" },
    @{ Find = "TERE"; Replace = "QWER" },
    @{ Find = "Tre";  Replace = "Qwer" },
    @{ Find = "Tre";  Replace = "Qwer" },
    @{ Find = "Tre";  Replace = "Qewr" },
    @{ Find = "Tre";  Replace = "Qewr" },
    @{ Find = "Tre";  Replace = "Qwer" },
    @{ Find = "tre";  Replace = "qwer" },
    @{ Find = "tre";  Replace = "qwer" },
    @{ Find = "tre";  Replace = "qwer" }
)

foreach ($item in $headerReplacements) {
    $rng.End = $hdrEnd
    $rng.Find.Execute($item.Find, $true, $true, $false, $false, $false, $true, 1, $false, $item.Replace, 1) | Out-Null
    $rng.Collapse(0)
}
